$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "_old" / "_new" AHB-version-suffixed header columns to the
#    concrete format-version labels used by this merged AHB diff
#    (FV2310 = previous version, FV2404 = current version).
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"
# K1 stays "diff"
$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# 2. Turn the populated range into an Excel table ("Table1") so the
#    regenerated AHB diff is filterable/sortable like the other merged files.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U59"), [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# 3. Freeze the header row so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
